$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 132, shifting existing rows 132-138 down to 133-139
$ws.Rows.Item(132).Insert()

# Fill the new row 132 with data
$ws.Cells.Item(132, 1).Value = 11
$ws.Cells.Item(132, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(132, 3).Value = "Bíobío"
$ws.Cells.Item(132, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(132, 4).Value = 44747
$ws.Cells.Item(132, 5).Value = 8
$ws.Cells.Item(132, 6).Value = 100112043
$ws.Cells.Item(132, 7).Value = "Pepino ensalada"
$ws.Cells.Item(132, 8).Value = "Sin especificar"
$ws.Cells.Item(132, 9).Value = "Primera"
$ws.Cells.Item(132, 10).Value = 190
$ws.Cells.Item(132, 11).Value = 16000
$ws.Cells.Item(132, 12).Value = 17000
$ws.Cells.Item(132, 13).Value = 16526
$ws.Cells.Item(132, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(132, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(132, 16).Value = 331
$ws.Cells.Item(132, 17).Value = 50
$ws.Cells.Item(132, 18).Value = "Hortaliza"
